# Auto-generated edit script: updates Leve profit/price figures across all sheets
# per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 348.25
$ws.Range("I31").Value = 348.25
$ws.Range("K31").Value = 1044.75
$ws.Range("M31").Value = -814.75
$ws.Range("H112").Value = 26317404
$ws.Range("J112").Value = 1837.1562
$ws.Range("L112").Value = 5511.4686
$ws.Range("N112").Value = -7727.4686
$ws.Range("H137").Value = 2704184.5
$ws.Range("I137").Value = 5556593.5
$ws.Range("J137").Value = 1902.2632
$ws.Range("K137").Value = 16669780.5
$ws.Range("L137").Value = 5706.7896
$ws.Range("M137").Value = -16667230.5
$ws.Range("N137").Value = -10806.7896
$ws.Range("H138").Value = 3010.875
$ws.Range("J138").Value = 2667.5789
$ws.Range("L138").Value = 8002.736699999999
$ws.Range("N138").Value = -18282.7367

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 950
$ws.Range("I19").Value = 950
$ws.Range("K19").Value = 950
$ws.Range("M19").Value = -721
$ws.Range("H61").Value = 41751350
$ws.Range("I61").Value = 50051216
$ws.Range("J61").Value = 252025
$ws.Range("K61").Value = 50051216
$ws.Range("L61").Value = 252025
$ws.Range("M61").Value = -50051004
$ws.Range("N61").Value = -252449
$ws.Range("H74").Value = 6462672.5
$ws.Range("I74").Value = 11410288
$ws.Range("J74").Value = 59875.293
$ws.Range("K74").Value = 11410288
$ws.Range("L74").Value = 59875.293
$ws.Range("M74").Value = -11409414
$ws.Range("N74").Value = -61623.293
$ws.Range("H77").Value = 6462672.5
$ws.Range("I77").Value = 11410288
$ws.Range("J77").Value = 59875.293
$ws.Range("K77").Value = 57051440
$ws.Range("L77").Value = 299376.465
$ws.Range("M77").Value = -57047072
$ws.Range("N77").Value = -308112.465
$ws.Range("H132").Value = 65340.395
$ws.Range("I132").Value = 73020.78999999999
$ws.Range("J132").Value = 59681.156
$ws.Range("K132").Value = 219062.37
$ws.Range("L132").Value = 179043.468
$ws.Range("M132").Value = -216532.37
$ws.Range("N132").Value = -184103.468
$ws.Range("H136").Value = 41751350
$ws.Range("I136").Value = 50051216
$ws.Range("J136").Value = 252025
$ws.Range("K136").Value = 150153648
$ws.Range("L136").Value = 756075
$ws.Range("M136").Value = -150151098
$ws.Range("N136").Value = -761175

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2080.9
$ws.Range("I134").Value = 3004.2222
$ws.Range("J134").Value = 1325.4546
$ws.Range("K134").Value = 9012.6666
$ws.Range("L134").Value = 3976.3638
$ws.Range("M134").Value = -6477.6666
$ws.Range("N134").Value = -9046.363799999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1253.0754
$ws.Range("I31").Value = 794.78125
$ws.Range("J31").Value = 1951.4286
$ws.Range("K31").Value = 794.78125
$ws.Range("L31").Value = 1951.4286
$ws.Range("M31").Value = -499.78125
$ws.Range("N31").Value = -2541.4286
$ws.Range("H32").Value = 1700
$ws.Range("I32").Value = 1700
$ws.Range("K32").Value = 1700
$ws.Range("M32").Value = -1384
$ws.Range("H34").Value = 1253.0754
$ws.Range("I34").Value = 794.78125
$ws.Range("J34").Value = 1951.4286
$ws.Range("K34").Value = 794.78125
$ws.Range("L34").Value = 1951.4286
$ws.Range("M34").Value = -592.78125
$ws.Range("N34").Value = -2355.4286
$ws.Range("H58").Value = 15153099
$ws.Range("I58").Value = 22223552
$ws.Range("J58").Value = 2129.1904
$ws.Range("K58").Value = 22223552
$ws.Range("L58").Value = 2129.1904
$ws.Range("M58").Value = -22223349
$ws.Range("N58").Value = -2535.1904
$ws.Range("H132").Value = 22692.043
$ws.Range("I132").Value = 1874.5
$ws.Range("J132").Value = 44414.695
$ws.Range("K132").Value = 5623.5
$ws.Range("L132").Value = 133244.085
$ws.Range("M132").Value = -3093.5
$ws.Range("N132").Value = -138304.085
$ws.Range("H134").Value = 26093.137
$ws.Range("I134").Value = 1681.125
$ws.Range("J134").Value = 91191.836
$ws.Range("K134").Value = 5043.375
$ws.Range("L134").Value = 273575.508
$ws.Range("M134").Value = -2508.375
$ws.Range("N134").Value = -278645.508
$ws.Range("H136").Value = 15153099
$ws.Range("I136").Value = 22223552
$ws.Range("J136").Value = 2129.1904
$ws.Range("K136").Value = 66670656
$ws.Range("L136").Value = 6387.5712
$ws.Range("M136").Value = -66668106
$ws.Range("N136").Value = -11487.5712

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H131").Value = 1188.2787
$ws.Range("J131").Value = 1292.4509
$ws.Range("L131").Value = 3877.3527
$ws.Range("N131").Value = -13957.3527

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1713.1
$ws.Range("I97").Value = 1802.5
$ws.Range("J97").Value = 1355.5
$ws.Range("K97").Value = 1802.5
$ws.Range("L97").Value = 1355.5
$ws.Range("M97").Value = -1306.5
$ws.Range("N97").Value = -2347.5
$ws.Range("H102").Value = 2083
$ws.Range("I102").Value = 2075.8572
$ws.Range("J102").Value = 2099.6667
$ws.Range("K102").Value = 2075.8572
$ws.Range("L102").Value = 2099.6667
$ws.Range("M102").Value = -453.8571999999999
$ws.Range("N102").Value = -5343.6667
$ws.Range("H132").Value = 183564.55
$ws.Range("I132").Value = 334004
$ws.Range("J132").Value = 127149.75
$ws.Range("K132").Value = 1002012
$ws.Range("L132").Value = 381449.25
$ws.Range("M132").Value = -999482
$ws.Range("N132").Value = -386509.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H132").Value = 21162.04
$ws.Range("I132").Value = 987.85
$ws.Range("J132").Value = 101858.8
$ws.Range("K132").Value = 2963.55
$ws.Range("L132").Value = 305576.4
$ws.Range("M132").Value = -433.5500000000002
$ws.Range("N132").Value = -310636.4
$ws.Range("H133").Value = 29973.889
$ws.Range("J133").Value = 29973.889
$ws.Range("L133").Value = 29973.889
$ws.Range("N133").Value = -35033.889
$ws.Range("H136").Value = 167919.83
$ws.Range("I136").Value = 101303.8
$ws.Range("J136").Value = 501000
$ws.Range("K136").Value = 303911.4
$ws.Range("L136").Value = 1503000
$ws.Range("M136").Value = -301361.4
$ws.Range("N136").Value = -1508100

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 84722.125
$ws.Range("I132").Value = 53501.74
$ws.Range("J132").Value = 203359.6
$ws.Range("K132").Value = 160505.22
$ws.Range("L132").Value = 610078.8
$ws.Range("M132").Value = -157975.22
$ws.Range("N132").Value = -615138.8
$ws.Range("H136").Value = 58242.23
$ws.Range("I136").Value = 30543.47
$ws.Range("J136").Value = 1000000
$ws.Range("K136").Value = 91630.41
$ws.Range("L136").Value = 3000000
$ws.Range("M136").Value = -89080.41
$ws.Range("N136").Value = -3005100

Write-Output "Applied 180 cell updates across 8 sheets."
